# Auto-generated Excel COM-interop script applying scheduled market-data refresh
# to the Leve profit tracker workbook (columns H-N per sheet/row, as produced by the
# scheduled runner's scrape + recompute pass).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 17618.715
$ws.Range("I6").Value = 100000
$ws.Range("J6").Value = 3888.5
$ws.Range("K6").Value = 300000
$ws.Range("L6").Value = 11665.5
$ws.Range("M6").Value = -299888
$ws.Range("N6").Value = -11889.5

$ws.Range("H9").Value = 195.5
$ws.Range("I9").Value = 47
$ws.Range("J9").Value = 245
$ws.Range("K9").Value = 47
$ws.Range("L9").Value = 245
$ws.Range("M9").Value = 122
$ws.Range("N9").Value = -583

$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("M12").ClearContents()
$ws.Range("N12").ClearContents()

$ws.Range("H17").Value = 3715.4
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 3715.4
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 11146.2
$ws.Range("N17").Value = -11482.2

$ws.Range("H70").Value = 5596.125
$ws.Range("I70").Value = 5000
$ws.Range("J70").Value = 5953.8
$ws.Range("K70").Value = 15000
$ws.Range("L70").Value = 17861.4
$ws.Range("M70").Value = -14730
$ws.Range("N70").Value = -18401.4

$ws.Range("H73").Value = 5596.125
$ws.Range("I73").Value = 5000
$ws.Range("J73").Value = 5953.8
$ws.Range("K73").Value = 15000
$ws.Range("L73").Value = 17861.4
$ws.Range("M73").Value = -14064
$ws.Range("N73").Value = -19733.4

$ws.Range("H93").Value = 40999.5
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 40999.5
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 40999.5
$ws.Range("N93").Value = -45991.5

$ws.Range("H132").Value = 5312.3887
$ws.Range("I132").Value = 1464.2174
$ws.Range("J132").Value = 12120.692
$ws.Range("K132").Value = 4392.6522
$ws.Range("L132").Value = 36362.076
$ws.Range("M132").Value = -1862.6522

$ws.Range("H137").Value = 325832.7
$ws.Range("I137").Value = 558561.5600000001
$ws.Range("J137").Value = 3592.6924
$ws.Range("K137").Value = 1675684.68
$ws.Range("L137").Value = 10778.0772
$ws.Range("M137").Value = -1673134.68

$ws.Range("H138").Value = 6044.2925
$ws.Range("I138").Value = 2102
$ws.Range("J138").Value = 6720.1143
$ws.Range("K138").Value = 6306
$ws.Range("L138").Value = 20160.3429
$ws.Range("M138").Value = -1166
$ws.Range("N138").Value = -30440.3429

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 964.2857
$ws.Range("I2").Value = 964.2857
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 964.2857
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -851.2857

$ws.Range("H74").Value = 15625930
$ws.Range("I74").Value = 19231544
$ws.Range("J74").Value = 1598.6666
$ws.Range("K74").Value = 19231544
$ws.Range("L74").Value = 1598.6666
$ws.Range("M74").Value = -19230670

$ws.Range("H77").Value = 15625930
$ws.Range("I77").Value = 19231544
$ws.Range("J77").Value = 1598.6666
$ws.Range("K77").Value = 96157720
$ws.Range("L77").Value = 7993.333000000001
$ws.Range("M77").Value = -96153352

$ws.Range("H92").Value = 20000000
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 20000000
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 20000000
$ws.Range("N92").Value = -20004992

$ws.Range("H116").Value = 964.2857
$ws.Range("I116").Value = 964.2857
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 964.2857
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = 1329.7143

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 964.2857
$ws.Range("I3").Value = 964.2857
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 964.2857
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -850.2857

$ws.Range("H20").Value = 8513.299999999999
$ws.Range("I20").Value = 10650
$ws.Range("J20").Value = 6765.091
$ws.Range("K20").Value = 10650
$ws.Range("L20").Value = 6765.091
$ws.Range("M20").Value = -10403
$ws.Range("N20").Value = -7259.091

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 270.93332
$ws.Range("I7").Value = 202.3
$ws.Range("J7").Value = 408.2
$ws.Range("K7").Value = 202.3
$ws.Range("L7").Value = 408.2
$ws.Range("M7").Value = -89.30000000000001
$ws.Range("N7").Value = -634.2

$ws.Range("H22").Value = 1632.4
$ws.Range("I22").Value = 1731.5
$ws.Range("J22").Value = 1566.3334
$ws.Range("K22").Value = 1731.5
$ws.Range("L22").Value = 1566.3334
$ws.Range("M22").Value = -1381.5
$ws.Range("N22").Value = -2266.3334

$ws.Range("H31").Value = 15387904
$ws.Range("I31").Value = 16669583
$ws.Range("J31").Value = 7753.8
$ws.Range("K31").Value = 16669583
$ws.Range("L31").Value = 7753.8
$ws.Range("M31").Value = -16669288

$ws.Range("H34").Value = 15387904
$ws.Range("I34").Value = 16669583
$ws.Range("J34").Value = 7753.8
$ws.Range("K34").Value = 16669583
$ws.Range("L34").Value = 7753.8
$ws.Range("M34").Value = -16669381

$ws.Range("H141").Value = 122441.23
$ws.Range("I141").Value = 83999.8
$ws.Range("J141").Value = 128848.13
$ws.Range("K141").Value = 83999.8
$ws.Range("L141").Value = 128848.13
$ws.Range("M141").Value = -78819.8
$ws.Range("N141").Value = -139208.13

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 55.125
$ws.Range("I2").Value = 57.25
$ws.Range("J2").Value = 53
$ws.Range("K2").Value = 343.5
$ws.Range("L2").Value = 318
$ws.Range("M2").Value = -230.5
$ws.Range("N2").Value = -544

$ws.Range("H33").Value = 330.5
$ws.Range("I33").Value = 448.75
$ws.Range("J33").Value = 94
$ws.Range("K33").Value = 2692.5
$ws.Range("L33").Value = 564
$ws.Range("M33").Value = -2409.5

$ws.Range("H113").Value = 718.2727
$ws.Range("I113").Value = 999
$ws.Range("J113").Value = 690.2
$ws.Range("K113").Value = 2997
$ws.Range("L113").Value = 2070.6
$ws.Range("M113").Value = -827
$ws.Range("N113").Value = -6410.6

$ws.Range("H121").Value = 905.875
$ws.Range("I121").Value = 296
$ws.Range("J121").Value = 1046.6154
$ws.Range("K121").Value = 888
$ws.Range("L121").Value = 3139.8462
$ws.Range("M121").Value = 422
$ws.Range("N121").Value = -5759.8462

$ws.Range("H126").Value = 29249
$ws.Range("I126").Value = 28500
$ws.Range("J126").Value = 29998
$ws.Range("K126").Value = 85500
$ws.Range("L126").Value = 89994
$ws.Range("M126").Value = -80560

$ws.Range("H128").Value = 139666
$ws.Range("I128").Value = 139666
$ws.Range("J128").Value = 0
$ws.Range("K128").Value = 418998
$ws.Range("L128").Value = 0
$ws.Range("M128").Value = -414018

$ws.Range("H132").Value = 1517.8649
$ws.Range("I132").Value = 1124.3914
$ws.Range("J132").Value = 2164.2856
$ws.Range("K132").Value = 10119.5226
$ws.Range("L132").Value = 19478.5704
$ws.Range("M132").Value = -7589.5226
$ws.Range("N132").Value = -24538.5704

$ws.Range("H139").Value = 2274.7307
$ws.Range("I139").Value = 1257.9445
$ws.Range("J139").Value = 4562.5
$ws.Range("K139").Value = 3773.8335
$ws.Range("L139").Value = 13687.5
$ws.Range("M139").Value = 1366.1665
$ws.Range("N139").Value = -23967.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 302.43243
$ws.Range("I2").Value = 243.47058
$ws.Range("J2").Value = 352.55
$ws.Range("K2").Value = 243.47058
$ws.Range("L2").Value = 352.55
$ws.Range("M2").Value = -130.47058
$ws.Range("N2").Value = -578.55

$ws.Range("H31").Value = 2630.25
$ws.Range("I31").Value = 2630.25
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 2630.25
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -2338.25

$ws.Range("H37").Value = 2630.25
$ws.Range("I37").Value = 2630.25
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 2630.25
$ws.Range("L37").Value = 0
$ws.Range("M37").Value = -2353.25

$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("M70").ClearContents()

$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("M73").ClearContents()

$ws.Range("H92").Value = 17136.092
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 17136.092
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 17136.092
$ws.Range("N92").Value = -20880.092

$ws.Range("H102").Value = 25004940
$ws.Range("I102").Value = 38465080
$ws.Range("J102").Value = 7539.143
$ws.Range("K102").Value = 38465080
$ws.Range("L102").Value = 7539.143
$ws.Range("M102").Value = -38463458
$ws.Range("N102").Value = -10783.143

$ws.Range("H126").Value = 13344.571
$ws.Range("I126").Value = 21256
$ws.Range("J126").Value = 10180
$ws.Range("K126").Value = 63768
$ws.Range("L126").Value = 30540
$ws.Range("M126").Value = -61298
$ws.Range("N126").Value = -35480

$ws.Range("H132").Value = 60604.53
$ws.Range("I132").Value = 68248.5
$ws.Range("J132").Value = 3274.75
$ws.Range("K132").Value = 204745.5
$ws.Range("L132").Value = 9824.25
$ws.Range("M132").Value = -202215.5
$ws.Range("N132").Value = -14884.25

$ws.Range("H134").Value = 58346.715
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 58346.715
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 175040.145
$ws.Range("N134").Value = -180110.145

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 916.8125
$ws.Range("I22").Value = 790
$ws.Range("J22").Value = 1466.3334
$ws.Range("K22").Value = 790
$ws.Range("L22").Value = 1466.3334
$ws.Range("M22").Value = -495
$ws.Range("N22").Value = -2056.3334

$ws.Range("H27").Value = 916.8125
$ws.Range("I27").Value = 790
$ws.Range("J27").Value = 1466.3334
$ws.Range("K27").Value = 790
$ws.Range("L27").Value = 1466.3334
$ws.Range("M27").Value = -683
$ws.Range("N27").Value = -1680.3334

$ws.Range("H61").Value = 4108.9165
$ws.Range("I61").Value = 4394.5557
$ws.Range("J61").Value = 3252
$ws.Range("K61").Value = 4394.5557
$ws.Range("L61").Value = 3252
$ws.Range("M61").Value = -4192.5557

$ws.Range("H113").Value = 4108.9165
$ws.Range("I113").Value = 4394.5557
$ws.Range("J113").Value = 3252
$ws.Range("K113").Value = 4394.5557
$ws.Range("L113").Value = 3252
$ws.Range("M113").Value = -2224.5557

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 13380
$ws.Range("I54").Value = 13380
$ws.Range("J54").Value = 0
$ws.Range("K54").Value = 13380
$ws.Range("L54").Value = 0
$ws.Range("M54").Value = -12860

$ws.Range("H96").Value = 2296.6667
$ws.Range("I96").Value = 2296.6667
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 2296.6667
$ws.Range("L96").Value = 0
$ws.Range("M96").Value = -923.6667000000002

$ws.Range("H100").Value = 4540.074
$ws.Range("I100").Value = 511.9375
$ws.Range("J100").Value = 10399.182
$ws.Range("K100").Value = 1023.875
$ws.Range("L100").Value = 20798.364
$ws.Range("M100").Value = -482.875

$ws.Range("H113").Value = 499.3913
$ws.Range("I113").Value = 437.0625
$ws.Range("J113").Value = 641.8570999999999
$ws.Range("K113").Value = 1311.1875
$ws.Range("L113").Value = 1925.5713
$ws.Range("M113").Value = 858.8125
$ws.Range("N113").Value = -6265.5713

$ws.Range("H126").Value = 40002468
$ws.Range("I126").Value = 50002308
$ws.Range("J126").Value = 3118
$ws.Range("K126").Value = 150006924
$ws.Range("L126").Value = 9354
$ws.Range("M126").Value = -150004454
